# Add new "data_collection_mode list" sheet (values DDA/DIA/MRM/PRM) positioned
# right after "ion_mobility list" and before "column_length_unit list", matching
# the new sheet order: ... ion_mobility list, data_collection_mode list,
# column_length_unit list, column_temp_unit list, spatial_type list,
# spatial_sampling_type list, resolution_x_unit list, resolution_y_unit list.

$wb = $excel.ActiveWorkbook

$ionMobilitySheet = $wb.Worksheets.Item("ion_mobility list")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ionMobilitySheet)
$newSheet.Name = "data_collection_mode list"

$newSheet.Range("A1").Value = "DDA"
$newSheet.Range("A2").Value = "DIA"
$newSheet.Range("A3").Value = "MRM"
$newSheet.Range("A4").Value = "PRM"

# Update the comment on the data_collection_mode column (X1) in the main sheet
# to mention the two new values.
$ws = $wb.Worksheets.Item("Export as TSV")
$commentCell = $ws.Range("X1")
$commentCell.Comment.Text("Mode of data collection in tandem MS assays. Either DDA (Data-dependent acquisition), DIA (Data-independent acquisition), MRM (multiple reaction monitoring), or PRM (parallel reaction monitoring).") | Out-Null

# Add data validation on column X (data_collection_mode) restricting to the new list.
$dataCollectionModeRange = $ws.Range("X2:X1048576")
$dataCollectionModeRange.Validation.Add(3, 1, 1, "='data_collection_mode list'!`$A`$1:`$A`$4")
$dataCollectionModeRange.Validation.ErrorTitle = "Value must come from list"
$dataCollectionModeRange.Validation.ErrorMessage = "Value must be one of: DDA / DIA / MRM / PRM."
$dataCollectionModeRange.Validation.ShowInput = $true
$dataCollectionModeRange.Validation.ShowError = $true

# Restore the originally-active tab ("Export as TSV") as the selected sheet,
# since adding a new worksheet makes it active by default.
$ws.Activate()
